{"js": "// Translate the English source strings in the document body (and the\n// review comment) to Spanish, matching the Crowdin \"es\" translation.\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText, options) {\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Language switcher line (top of the document).\nawait replaceOnce(\"English\", \"Ingl\u00e9s\");\nawait replaceOnce(\n  \" / Portuguese / French / Thai / Vietnamese / Spanish\",\n  \" / Portugu\u00e9s / Franc\u00e9s / Tailand\u00e9s / Vietnamita / Espa\u00f1ol\"\n);\n\n// \"English\" heading that labels this language section.\nawait replaceOnce(\"English\", \"Ingl\u00e9s\");\n\n// Brief / target audience table.\nawait replaceOnce(\"Brief\", \"Breve\");\nawait replaceOnce(\n  \"An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io\",\n  \"An email sent to partners in the target country whose documents failed our verification process. Se enviar\u00e1 a trav\u00e9s de customer.io\"\n);\nawait replaceOnce(\"Target audience\", \"P\u00fablico objetivo\");\n\n// Email body copy.\nawait replaceOnce(\n  \"Uh oh! Your documents couldn\\u2019t be verified\",\n  \"\\u00a1Oh no! Tus documentos no han podido ser verificados\"\n);\nawait replaceOnce(\"Hi \", \"Hola \");\nawait replaceOnce(\n  \"We regret to inform you that your documents have failed our verification process as we found the following issues with them: \",\n  \"Lamentamos informarte de que tus documentos no han superado nuestro proceso de verificaci\u00f3n, ya que hemos encontrado los siguientes problemas en ellos: \"\n);\nawait replaceOnce(\n  \"A copy of your vaccination certificate\",\n  \"Una copia de tu certificado de vacunaci\u00f3n\"\n);\nawait replaceOnce(\": Document is unclear\", \": El documento no est\u00e1 claro\");\nawait replaceOnce(\"[Document 2]\", \"[Documento 2]\");\nawait replaceOnce(\n  \"Please resubmit the documents above by \",\n  \"Por favor, vuelve a enviar los documentos mencionados arriba antes del \"\n);\nawait replaceOnce(\n  \" so we can proceed with the necessary arrangements.\",\n  \" para que podamos proceder con los preparativos necesarios.\"\n);\nawait replaceOnce(\n  \"If you have any questions, please contact us via \",\n  \"Si tienes alguna pregunta, entra en contacto con nosotros por \"\n);\nawait replaceOnce(\" or \", \" o \"); // live chat / WhatsApp sentence\nawait replaceOnce(\n  \"If you have any questions, please contact your country manager, \",\n  \"Si tienes alguna pregunta, entra en contacto con el gestor de tu pa\u00eds \"\n);\nawait replaceOnce(\", at \", \", en \");\nawait replaceOnce(\" or \", \" o \"); // [EMAIL ADDRESS] / [WHATSAPP NO] sentence\n\n// Review comment text.\nconst comments = context.document.comments;\ncomments.load(\"items\");\nawait context.sync();\nfor (const comment of comments.items) {\n  comment.load(\"content\");\n}\nawait context.sync();\nfor (const comment of comments.items) {\n  if (comment.content === \"choose either one\") {\n    comment.content = \"elija uno de los dos\";\n  }\n}\nawait context.sync();\n", "ps1": "# Translate the English source strings in the document body (and the\n# review comment) to Spanish, matching the Crowdin \"es\" translation.\n\n$d = $word.ActiveDocument\n\n# Word Find/Replace constants.\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\nfunction Replace-FirstMatch($findText, $replaceText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        $wdReplaceOne # Replace\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# Language switcher line (top of the document).\nReplace-FirstMatch \"English\" \"Ingl\u00e9s\"\nReplace-FirstMatch \" / Portuguese / French / Thai / Vietnamese / Spanish\" \" / Portugu\u00e9s / Franc\u00e9s / Tailand\u00e9s / Vietnamita / Espa\u00f1ol\"\n\n# \"English\" heading that labels this language section.\nReplace-FirstMatch \"English\" \"Ingl\u00e9s\"\n\n# Brief / target audience table.\nReplace-FirstMatch \"Brief\" \"Breve\"\nReplace-FirstMatch \"An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io\" \"An email sent to partners in the target country whose documents failed our verification process. Se enviar\u00e1 a trav\u00e9s de customer.io\"\nReplace-FirstMatch \"Target audience\" \"P\u00fablico objetivo\"\n\n# Email body copy.\nReplace-FirstMatch \"Uh oh! Your documents couldn\u2019t be verified\" \"\u00a1Oh no! Tus documentos no han podido ser verificados\"\nReplace-FirstMatch \"Hi \" \"Hola \"\nReplace-FirstMatch \"We regret to inform you that your documents have failed our verification process as we found the following issues with them: \" \"Lamentamos informarte de que tus documentos no han superado nuestro proceso de verificaci\u00f3n, ya que hemos encontrado los siguientes problemas en ellos: \"\nReplace-FirstMatch \"A copy of your vaccination certificate\" \"Una copia de tu certificado de vacunaci\u00f3n\"\nReplace-FirstMatch \": Document is unclear\" \": El documento no est\u00e1 claro\"\nReplace-FirstMatch \"[Document 2]\" \"[Documento 2]\"\nReplace-FirstMatch \"Please resubmit the documents above by \" \"Por favor, vuelve a enviar los documentos mencionados arriba antes del \"\nReplace-FirstMatch \" so we can proceed with the necessary arrangements.\" \" para que podamos proceder con los preparativos necesarios.\"\nReplace-FirstMatch \"If you have any questions, please contact us via \" \"Si tienes alguna pregunta, entra en contacto con nosotros por \"\nReplace-FirstMatch \" or \" \" o \"   # live chat / WhatsApp sentence\nReplace-FirstMatch \"If you have any questions, please contact your country manager, \" \"Si tienes alguna pregunta, entra en contacto con el gestor de tu pa\u00eds \"\nReplace-FirstMatch \", at \" \", en \"\nReplace-FirstMatch \" or \" \" o \"   # [EMAIL ADDRESS] / [WHATSAPP NO] sentence\n\n# Review comment text.\nforeach ($c in $d.Comments) {\n    if ($c.Range.Text -eq \"choose either one\") {\n        $c.Range.Text = \"elija uno de los dos\"\n    }\n}\n"}
